$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C48").Value = 104078
$ws.Range("E48").Value = 624994490

$ws.Range("C65").Value = 61058
$ws.Range("E65").Value = 334115122

$ws.Range("C81").Value = 26166
$ws.Range("E81").Value = 165097918

$ws.Range("C99").Value = 136576
$ws.Range("E99").Value = 863132204

$ws.Range("C160").Value = 26949
$ws.Range("D160").Value = 11393
$ws.Range("E160").Value = 37313966

$ws.Range("C164").Value = 50579
$ws.Range("E164").Value = 168930017

$ws.Range("C168").Value = 285052
$ws.Range("E168").Value = 1211795750

$ws.Range("C169").Value = 562628
$ws.Range("E169").Value = 1285106339

$ws.Range("C170").Value = 367457
$ws.Range("E170").Value = 2846922007

$ws.Range("C171").Value = 115187
$ws.Range("E171").Value = 447452079

$ws.Range("C174").Value = 357282
$ws.Range("D174").Value = 69790
$ws.Range("E174").Value = 1018803967

$ws.Range("C178").Value = 75361
$ws.Range("E178").Value = 102748863

$ws.Range("C179").Value = 235741
$ws.Range("E179").Value = 813061344

$ws.Range("C180").Value = 141503
$ws.Range("E180").Value = 341108088

$ws.Range("C262").Value = 38984
$ws.Range("E262").Value = 124727562

$ws.Range("C264").Value = 47472
$ws.Range("E264").Value = 81956515

$ws.Range("C279").Value = 28965
$ws.Range("E279").Value = 57087380

$ws.Range("C312").Value = 75102
$ws.Range("E312").Value = 201349936

$ws.Range("C320").Value = 67243
$ws.Range("E320").Value = 124557315
